$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 25: 2025-12-02 (serial 45993), Error Count 86
$ws.Range("A24").Copy($ws.Range("A25"))
$ws.Range("A25").Value = 45993
$ws.Range("B25").Value = 86

# New row 26: 2025-12-03 (serial 45994), Error Count 85
$ws.Range("A24").Copy($ws.Range("A26"))
$ws.Range("A26").Value = 45994
$ws.Range("B26").Value = 85

$ws.Range("A25:B26").Select() | Out-Null
